$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain text (source data is inline strings,
# some of which (e.g. "29.890.54") are not valid numbers) while we assign values,
# then restore the original (default) cell style so no formatting is left behind.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range('D2').Value = '29.890.54'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '1.632.20'
$ws.Range('E3').Value = '  +1.56%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '214.58'
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '28.64'
$ws.Range('E8').Value = '  +1.88%  '
$ws.Range('E9').Value = '  +2.08%  '
$ws.Range('E10').Value = '  +0.77%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '1.865.30'
$ws.Range('E12').Value = '  +1.51%  '
$ws.Range('D13').Value = '1.630.49'
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('D15').Value = '9.30'
$ws.Range('E15').Value = '  +18.44%  '
$ws.Range('E16').Value = '  +2.60%  '
$ws.Range('D17').Value = '29.887.95'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '242.79'
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = '9.85'
$ws.Range('E22').Value = '  +4.64%  '
$ws.Range('E23').Value = '  +2.19%  '
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('D25').Value = '157.45'
$ws.Range('D26').Value = '15.53'
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('D28').Value = '6.59'
$ws.Range('E28').Value = '  +1.86%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('E31').Value = '  +4.19%  '
$ws.Range('E32').Value = '  +3.71%  '
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('D34').Value = '1.429.36'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('E35').Value = '  +4.66%  '
$ws.Range('D36').Value = '1.04'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = '2.80'
$ws.Range('E37').Value = '  -4.61%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('D40').Value = '75.58'
$ws.Range('E40').Value = '  +13.88%  '
$ws.Range('D41').Value = '0.551'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('E42').Value = '  +1.96%  '
$ws.Range('D43').Value = '0.825'
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('D45').Value = '53.65'
$ws.Range('E45').Value = '  -5.66%  '
$ws.Range('E46').Value = '  +3.61%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = '5.37'
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('D49').Value = '1.773.28'
$ws.Range('E49').Value = '  +1.75%  '
$ws.Range('D50').Value = '89.04'
$ws.Range('E50').Value = '  +2.75%  '
$ws.Range('D51').Value = '0.0₆0111'
$ws.Range('E51').Value = '  +5.82%  '

$rng.Style = "Normal"
